# Updates the cryptos list: refreshed prices and Volume(1h) percentages,
# plus coin-rank swaps in rows 22/23, 25/26, and 47/48/49.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.230.12"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").Value = "3.926.23"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "486.71"
$ws.Range("E5").Value = "  +3.35%  "
$ws.Range("D6").Value = "147.58"
$ws.Range("E6").Value = "  +1.43%  "
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D9").Value = "0.735"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("E10").Value = "  +1.93%  "
$ws.Range("E11").Value = "  +4.84%  "
$ws.Range("D12").Value = "43.02"
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").Value = "10.73"
$ws.Range("E13").Value = "  +3.32%  "
$ws.Range("D14").Value = "4.546.20"
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").Value = "3.928.59"
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").Value = "14.42"
$ws.Range("E16").Value = "  -4.87%  "
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "19.94"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("E19").Value = "  -1.99%  "
$ws.Range("D20").Value = "68.379.53"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").Value = "442.81"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "15.16"
$ws.Range("E22").Value = "  +4.40%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").Value = "3.49"
$ws.Range("E23").Value = "  +3.03%  "
$ws.Range("D24").Value = "88.34"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").Value = "11.31"
$ws.Range("E25").Value = "  +15.17%  "
$ws.Range("B26").Value = "Filecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D26").Value = "11.47"
$ws.Range("E26").Value = "  +11.45%  "
$ws.Range("D27").Value = "3.63"
$ws.Range("E27").Value = "  +0.90%  "
$ws.Range("D28").Value = "38.82"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "5.73"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").Value = "719.43"
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("D31").Value = "13.77"
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("D32").Value = "0.131"
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("D33").Value = "2.92"
$ws.Range("E33").Value = "  +3.73%  "
$ws.Range("D34").Value = "6.42"
$ws.Range("E34").Value = "  +19.73%  "
$ws.Range("D35").Value = "42.58"
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("D36").Value = "0.0₃0891"
$ws.Range("E36").Value = "  +13.83%  "
$ws.Range("D37").Value = "61.04"
$ws.Range("E37").Value = "  +5.46%  "
$ws.Range("D38").Value = "0.150"
$ws.Range("E38").Value = "  -1.81%  "
$ws.Range("D39").Value = "0.403"
$ws.Range("E39").Value = "  +20.00%  "
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("D41").Value = "3.02"
$ws.Range("E41").Value = "  +17.25%  "
$ws.Range("D42").Value = "3.26"
$ws.Range("E42").Value = "  +7.61%  "
$ws.Range("D43").Value = "0.0482"
$ws.Range("E43").Value = "  +1.24%  "
$ws.Range("D44").Value = "2.90"
$ws.Range("E44").Value = "  +2.95%  "
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "3.28"
$ws.Range("E47").Value = "  +2.94%  "
$ws.Range("B48").Value = "LidoDAOToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D48").Value = "3.42"
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0351"
$ws.Range("E49").Value = "  +32.33%  "
$ws.Range("D50").Value = "2.15"
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("D51").Value = "145.85"
$ws.Range("E51").Value = "  -0.70%  "
